$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (A-level name / checkpoint text, B, C, D=computed diff)
$rows = @(
    @{ R = 136; A = "1st Move";          B = 45752; C = 53228 },
    @{ R = 137; A = "Checkpoint 89";     B = 45830; C = 53306 },
    @{ R = 138; A = "Checkpoint 404";    B = 45929; C = 53405 },
    @{ R = 139; A = "Checkpoint 1037";   B = 46124; C = 53602 },
    @{ R = 140; A = "Checkpoint 1534";   B = 46279; C = 53758 },
    @{ R = 141; A = "Checkpoint 1836";   B = 46376; C = 53854 },
    @{ R = 142; A = "Checkpoitn 2224";   B = 46504; C = 53982 },
    @{ R = 143; A = "Checkpoint 2586";   B = 46624; C = 54102 },
    @{ R = 144; A = "Enter door";        B = 46876; C = 54354 },
    @{ R = 145; A = "Touch button";      B = 47892; C = 55387 },
    @{ R = 146; A = "End level";         B = 48805; C = 56300 },
    @{ R = 147; A = "Enter 8-5";         B = 50431; C = 59573 },
    @{ R = 148; A = "1st Move";          B = 50658; C = 59821 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
}

# Fill the D column with the same IF formula as the existing rows (D90:D135),
# continuing the fill-down pattern as one shared-formula range D136:D148.
$ws.Range("D136:D148").Formula = "=IF(B136 >  0,C136-B136, 0)"

$ws.Range("B149").Select() | Out-Null

Write-Host "done"
